# Update cryptos list values (prices in column D, 1h volume/change % in column E)
# per the commit "Updated cryptos list on Sat Jan 20 07:41:51 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.534.29"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "2.484.49"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.47"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.35"
$ws.Range("E6").Value = "  -2.40%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  +3.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.82"
$ws.Range("E10").Value = "  -2.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.111"
$ws.Range("E12").Value = "  +2.10%  "

$ws.Range("D13").Value = "2.865.83"
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.94"
$ws.Range("E14").Value = "  -1.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.43"
$ws.Range("E15").Value = "  +9.60%  "

$ws.Range("D16").Value = "2.474.19"
$ws.Range("E16").Value = "  -2.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").Value = "41.567.25"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.56"
$ws.Range("E19").Value = "  +3.76%  "

$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("E20").Value = "  +2.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.52"
$ws.Range("E21").Value = "  +5.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.19"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.92"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.73"
$ws.Range("E24").Value = "  -0.68%  "

$ws.Range("E25").Value = "  +0.55%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.07"
$ws.Range("E27").Value = "  +4.10%  "

$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("E29").Value = "  +0.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.10"
$ws.Range("E30").Value = "  -1.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.50"
$ws.Range("E31").Value = "  +3.79%  "

$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0759"
$ws.Range("E34").Value = "  +1.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.52"
$ws.Range("E35").Value = "  +2.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.40"
$ws.Range("E36").Value = "  -9.71%  "

$ws.Range("E37").Value = "  +3.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.91"
$ws.Range("E38").Value = "  -4.61%  "

$ws.Range("E39").Value = "  -1.70%  "

$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.06"
$ws.Range("E41").Value = "  -3.81%  "

$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("D43").Value = "1.973.19"
$ws.Range("E43").Value = "  -0.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.14"
$ws.Range("E44").Value = "  -3.10%  "

$ws.Range("E45").Value = "  -0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("E46").Value = "  -2.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.97"
$ws.Range("E47").Value = "  +2.94%  "

$ws.Range("D48").Value = "2.721.98"
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "98.25"
$ws.Range("E49").Value = "  +1.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.36"
$ws.Range("E50").Value = "  -1.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.174"
$ws.Range("E51").Value = "  -2.31%  "
